$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 331.42856
$ws.Range("I2").Value = 424
$ws.Range("K2").Value = 424
$ws.Range("M2").Value = -311

$ws.Range("H18").Value = 1857.2727
$ws.Range("I18").Value = 1857.2727
$ws.Range("K18").Value = 1857.2727
$ws.Range("M18").Value = -1573.2727

$ws.Range("H32").Value = 4778.222
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4778.222
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4778.222
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -5430.222

$ws.Range("H41").Value = 536.1667
$ws.Range("I41").Value = 760.3333
$ws.Range("J41").Value = 461.44446
$ws.Range("K41").Value = 760.3333
$ws.Range("L41").Value = 461.44446
$ws.Range("M41").Value = -320.3333
$ws.Range("N41").Value = -1341.44446

$ws.Range("H76").Value = 3541.8948
$ws.Range("I76").Value = 3469
$ws.Range("J76").Value = 3622.889
$ws.Range("K76").Value = 3469
$ws.Range("L76").Value = 3622.889
$ws.Range("M76").Value = -3154
$ws.Range("N76").Value = -4252.889

$ws.Range("H79").Value = 3541.8948
$ws.Range("I79").Value = 3469
$ws.Range("J79").Value = 3622.889
$ws.Range("K79").Value = 3469
$ws.Range("L79").Value = 3622.889
$ws.Range("M79").Value = -2377
$ws.Range("N79").Value = -5806.889

$ws.Range("H98").Value = 1117.7142
$ws.Range("I98").Value = 904
$ws.Range("J98").Value = 2400
$ws.Range("K98").Value = 904
$ws.Range("L98").Value = 2400
$ws.Range("M98").Value = 594
$ws.Range("N98").Value = -5396

$ws.Range("H107").Value = 383.94736
$ws.Range("I107").Value = 235
$ws.Range("J107").Value = 801
$ws.Range("K107").Value = 235
$ws.Range("L107").Value = 801
$ws.Range("M107").Value = 1685
$ws.Range("N107").Value = -4641

$ws.Range("H116").Value = 1742.8572
$ws.Range("I116").Value = 1680
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 1680
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = 1762
$ws.Range("N116").Value = -8784

$ws.Range("H122").Value = 1117.7142
$ws.Range("I122").Value = 904
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 2712
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -262
$ws.Range("N122").Value = -12100

$ws.Range("H137").Value = 4763988.5
$ws.Range("I137").Value = 7144639.5
$ws.Range("J137").Value = 2686.1428
$ws.Range("K137").Value = 21433918.5
$ws.Range("L137").Value = 8058.428400000001
$ws.Range("M137").Value = -21431368.5
$ws.Range("N137").Value = -13158.4284

$ws.Range("H138").Value = 850871.4399999999
$ws.Range("I138").Value = 1341.3077
$ws.Range("J138").Value = 1464421
$ws.Range("K138").Value = 4023.9231
$ws.Range("L138").Value = 4393263
$ws.Range("M138").Value = 1116.0769
$ws.Range("N138").Value = -4403543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 70150.13
$ws.Range("I132").Value = 54545.05
$ws.Range("J132").Value = 94858.164
$ws.Range("K132").Value = 163635.15
$ws.Range("L132").Value = 284574.492
$ws.Range("M132").Value = -161105.15
$ws.Range("N132").Value = -289634.492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1287.875
$ws.Range("I99").Value = 738.8889
$ws.Range("J99").Value = 1993.7142
$ws.Range("K99").Value = 738.8889
$ws.Range("L99").Value = 1993.7142
$ws.Range("M99").Value = 759.1111
$ws.Range("N99").Value = -4989.7142

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H134").Value = 3333
$ws.Range("I134").Value = 3333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7464
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9750
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9750
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9750
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -10340

$ws.Range("H34").Value = 9750
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9750
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9750
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -10154

$ws.Range("H107").Value = 701.4375
$ws.Range("I107").Value = 256.77777
$ws.Range("J107").Value = 1273.1428
$ws.Range("K107").Value = 256.77777
$ws.Range("L107").Value = 1273.1428
$ws.Range("M107").Value = 1663.22223
$ws.Range("N107").Value = -5113.1428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1104
$ws.Range("I131").Value = 685
$ws.Range("J131").Value = 1383.3334
$ws.Range("K131").Value = 2055
$ws.Range("L131").Value = 4150.0002
$ws.Range("M131").Value = 2985
$ws.Range("N131").Value = -14230.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2625.5833
$ws.Range("I16").Value = 1137
$ws.Range("J16").Value = 19000
$ws.Range("K16").Value = 1137
$ws.Range("L16").Value = 19000
$ws.Range("M16").Value = -967
$ws.Range("N16").Value = -19340

$ws.Range("H22").Value = 664.7059
$ws.Range("I22").Value = 533.3333
$ws.Range("J22").Value = 812.5
$ws.Range("K22").Value = 533.3333
$ws.Range("L22").Value = 812.5
$ws.Range("M22").Value = -238.3333
$ws.Range("N22").Value = -1402.5

$ws.Range("H27").Value = 664.7059
$ws.Range("I27").Value = 533.3333
$ws.Range("J27").Value = 812.5
$ws.Range("K27").Value = 533.3333
$ws.Range("L27").Value = 812.5
$ws.Range("M27").Value = -426.3333
$ws.Range("N27").Value = -1026.5

$ws.Range("H46").Value = 966.6667
$ws.Range("I46").Value = 950
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 950
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -762
$ws.Range("N46").Value = -1376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1311.3077
$ws.Range("I126").Value = 1063
$ws.Range("J126").Value = 1870
$ws.Range("K126").Value = 3189
$ws.Range("L126").Value = 5610
$ws.Range("M126").Value = -719
$ws.Range("N126").Value = -10550
